# Append a new row (row 95) to each of the 4 worksheets, duplicating the
# last existing row (row 94) and only updating the "time" (column A) value.
$wb = $excel.ActiveWorkbook

for ($i = 1; $i -le $wb.Worksheets.Count; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # Duplicate row 94 (values + formatting) into the new row 95.
    $srcRow = $ws.Range("A94:I94")
    $dstRow = $ws.Range("A95:I95")
    $srcRow.Copy($dstRow)

    # Update the timestamp in column A for the newly added row.
    $ws.Range("A95").Value = 45881.43618055555
}
